# Auto-generated edit script to apply Phantom_Profits numeric corrections
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H40").Value = 2582
$ws.Range("I40").Value = 1374.5
$ws.Range("K40").Value = 1374.5
$ws.Range("M40").Value = -1199.5
$ws.Range("H74").Value = 3638.5557
$ws.Range("I74").Value = 3638.5557
$ws.Range("K74").Value = 3638.5557
$ws.Range("M74").Value = -2702.5557
$ws.Range("H77").Value = 3638.5557
$ws.Range("I77").Value = 3638.5557
$ws.Range("K77").Value = 18192.7785
$ws.Range("M77").Value = -13512.7785
$ws.Range("H80").Value = 3115.4443
$ws.Range("I80").Value = 2573
$ws.Range("J80").Value = 3549.4
$ws.Range("K80").Value = 7719
$ws.Range("L80").Value = 10648.2
$ws.Range("M80").Value = -6721
$ws.Range("N80").Value = -12644.2
$ws.Range("H83").Value = 3115.4443
$ws.Range("I83").Value = 2573
$ws.Range("J83").Value = 3549.4
$ws.Range("K83").Value = 23157
$ws.Range("L83").Value = 31944.6
$ws.Range("M83").Value = -18165
$ws.Range("N83").Value = -41928.60000000001
$ws.Range("H98").Value = 622.9
$ws.Range("I98").Value = 682.1111
$ws.Range("J98").Value = 90
$ws.Range("K98").Value = 682.1111
$ws.Range("L98").Value = 90
$ws.Range("M98").Value = 815.8889
$ws.Range("N98").Value = -3086
$ws.Range("H122").Value = 622.9
$ws.Range("I122").Value = 682.1111
$ws.Range("J122").Value = 90
$ws.Range("K122").Value = 2046.3333
$ws.Range("L122").Value = 270
$ws.Range("M122").Value = 403.6667000000002
$ws.Range("N122").Value = -5170
$ws.Range("H137").Value = 1649
$ws.Range("I137").Value = 1289.7307
$ws.Range("K137").Value = 3869.1921
$ws.Range("M137").Value = -1319.1921

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H13").Value = 5020000
$ws.Range("J13").Value = 40000
$ws.Range("L13").Value = 40000
$ws.Range("N13").Value = -40288
$ws.Range("H61").Value = 2945.3635
$ws.Range("J61").Value = 2800
$ws.Range("L61").Value = 2800
$ws.Range("N61").Value = -3224
$ws.Range("H74").Value = 1399.8334
$ws.Range("I74").Value = 1000.25
$ws.Range("J74").Value = 2199
$ws.Range("K74").Value = 1000.25
$ws.Range("L74").Value = 2199
$ws.Range("M74").Value = -126.25
$ws.Range("N74").Value = -3947
$ws.Range("H77").Value = 1399.8334
$ws.Range("I77").Value = 1000.25
$ws.Range("J77").Value = 2199
$ws.Range("K77").Value = 5001.25
$ws.Range("L77").Value = 10995
$ws.Range("M77").Value = -633.25
$ws.Range("N77").Value = -19731
$ws.Range("H92").Value = 49412.25
$ws.Range("J92").Value = 49412.25
$ws.Range("L92").Value = 49412.25
$ws.Range("N92").Value = -54404.25
$ws.Range("H122").Value = 1619.6154
$ws.Range("I122").Value = 1619.6154
$ws.Range("K122").Value = 4858.8462
$ws.Range("M122").Value = -2408.8462
$ws.Range("H136").Value = 2945.3635
$ws.Range("J136").Value = 2800
$ws.Range("L136").Value = 8400
$ws.Range("N136").Value = -13500

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H20").Value = 2692.111
$ws.Range("I20").Value = 1676.4286
$ws.Range("K20").Value = 1676.4286
$ws.Range("M20").Value = -1429.4286
$ws.Range("H86").Value = 3839.0715
$ws.Range("I86").Value = 3839.0715
$ws.Range("K86").Value = 3839.0715
$ws.Range("M86").Value = -2716.0715
$ws.Range("H89").Value = 3839.0715
$ws.Range("I89").Value = 3839.0715
$ws.Range("K89").Value = 19195.3575
$ws.Range("M89").Value = -13579.3575

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H86").Value = 7897.8
$ws.Range("I86").Value = 7489
$ws.Range("J86").Value = 8000
$ws.Range("K86").Value = 7489
$ws.Range("L86").Value = 8000
$ws.Range("M86").Value = -6366
$ws.Range("N86").Value = -10246
$ws.Range("H89").Value = 7897.8
$ws.Range("I89").Value = 7489
$ws.Range("J89").Value = 8000
$ws.Range("K89").Value = 37445
$ws.Range("L89").Value = 40000
$ws.Range("M89").Value = -31829
$ws.Range("N89").Value = -51232
$ws.Range("H99").Value = 3180.4546
$ws.Range("I99").Value = 2999.8333
$ws.Range("J99").Value = 3397.2
$ws.Range("K99").Value = 2999.8333
$ws.Range("L99").Value = 3397.2
$ws.Range("M99").Value = -1501.8333
$ws.Range("N99").Value = -6393.2
$ws.Range("H122").Value = 4385.375
$ws.Range("I122").Value = 4912.75
$ws.Range("J122").Value = 2803.25
$ws.Range("K122").Value = 14738.25
$ws.Range("L122").Value = 8409.75
$ws.Range("M122").Value = -12288.25
$ws.Range("N122").Value = -13309.75
$ws.Range("H126").Value = 3180.4546
$ws.Range("I126").Value = 2999.8333
$ws.Range("J126").Value = 3397.2
$ws.Range("K126").Value = 8999.499899999999
$ws.Range("L126").Value = 10191.6
$ws.Range("M126").Value = -6529.499899999999
$ws.Range("N126").Value = -15131.6

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H10").Value = 176.66667
$ws.Range("I10").Value = 176.66667
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 530.00001
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -391.00001
$ws.Range("N10").ClearContents()
$ws.Range("H68").Value = 648.1667
$ws.Range("I68").Value = 598
$ws.Range("J68").Value = 673.25
$ws.Range("K68").Value = 1794
$ws.Range("L68").Value = 2019.75
$ws.Range("M68").Value = -983
$ws.Range("N68").Value = -3641.75
$ws.Range("H71").Value = 648.1667
$ws.Range("I71").Value = 598
$ws.Range("J71").Value = 673.25
$ws.Range("K71").Value = 5382
$ws.Range("L71").Value = 6059.25
$ws.Range("M71").Value = -1326
$ws.Range("N71").Value = -14171.25
$ws.Range("H125").Value = 16683
$ws.Range("I125").Value = 16683
$ws.Range("K125").Value = 50049
$ws.Range("M125").Value = -45129
$ws.Range("H132").Value = 1252.2941
$ws.Range("I132").Value = 1224.3572
$ws.Range("J132").Value = 1382.6666
$ws.Range("K132").Value = 11019.2148
$ws.Range("L132").Value = 12443.9994
$ws.Range("M132").Value = -8489.2148
$ws.Range("N132").Value = -17503.9994

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H70").Value = 6050.6924
$ws.Range("I70").Value = 6050.6924
$ws.Range("K70").Value = 6050.6924
$ws.Range("M70").Value = -5780.6924
$ws.Range("H73").Value = 6050.6924
$ws.Range("I73").Value = 6050.6924
$ws.Range("K73").Value = 6050.6924
$ws.Range("M73").Value = -5114.6924
$ws.Range("H102").Value = 2339.4
$ws.Range("I102").Value = 2339.4
$ws.Range("K102").Value = 2339.4
$ws.Range("M102").Value = -717.4000000000001
$ws.Range("H122").Value = 2080.238
$ws.Range("I122").Value = 1984
$ws.Range("J122").Value = 2994.5
$ws.Range("K122").Value = 5952
$ws.Range("L122").Value = 8983.5
$ws.Range("M122").Value = -3502
$ws.Range("N122").Value = -13883.5

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H7").Value = 6267.1113
$ws.Range("I7").Value = 6261
$ws.Range("K7").Value = 6261
$ws.Range("M7").Value = -6149
$ws.Range("H22").Value = 882.5714
$ws.Range("I22").Value = 366.66666
$ws.Range("J22").Value = 1269.5
$ws.Range("K22").Value = 366.66666
$ws.Range("L22").Value = 1269.5
$ws.Range("M22").Value = -71.66665999999998
$ws.Range("N22").Value = -1859.5
$ws.Range("H27").Value = 882.5714
$ws.Range("I27").Value = 366.66666
$ws.Range("J27").Value = 1269.5
$ws.Range("K27").Value = 366.66666
$ws.Range("L27").Value = 1269.5
$ws.Range("M27").Value = -259.66666
$ws.Range("N27").Value = -1483.5
$ws.Range("H40").Value = 3987.75
$ws.Range("I40").Value = 3987.75
$ws.Range("K40").Value = 3987.75
$ws.Range("M40").Value = -3851.75
$ws.Range("H105").Value = 35615
$ws.Range("J105").Value = 35615
$ws.Range("L105").Value = 35615
$ws.Range("N105").Value = -42603
$ws.Range("H122").Value = 3888.2856
$ws.Range("I122").Value = 3483.8
$ws.Range("K122").Value = 10451.4
$ws.Range("M122").Value = -8001.400000000001
$ws.Range("H126").Value = 6267.1113
$ws.Range("I126").Value = 6261
$ws.Range("K126").Value = 18783
$ws.Range("M126").Value = -16313

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H132").Value = 5004.6875
$ws.Range("I132").Value = 5006.25
$ws.Range("K132").Value = 15018.75
$ws.Range("M132").Value = -12488.75
$ws.Range("H133").Value = 56000
$ws.Range("I133").Value = 50000
$ws.Range("J133").Value = 60000
$ws.Range("K133").Value = 50000
$ws.Range("L133").Value = 60000
$ws.Range("M133").Value = -44940
$ws.Range("N133").Value = -70120
